$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: new NPC "cleaner" (base) ---------------------------------
$ws.Range("A17").Value = 5200
$ws.Range("B17").Value = "cleaner_enhancer"
$ws.Range("B18").Value = "cleaner_R_enhancer"
$ws.Range("C17").Value = "클리너 인핸서"
$ws.Range("C18").Value = "깔끔한 클리너 인핸서"
$ws.Range("D17").Value = "enhancer"
$ws.Range("E17").Value = "Weakness"
$ws.Range("F17").Value = "SmithsSmithy"
$ws.Range("G17").Value = "doncina_skill_01"
$ws.Range("K17").Value = 1

# --- Row 18: new NPC "cleaner_R" (variant) -----------------------------
$ws.Range("A18").Value = 5201
$ws.Range("D18").Value = "enhancer"
$ws.Range("E18").Value = "Weakness"
$ws.Range("F18").Value = "SmithsSmithy"
$ws.Range("H18").Value = "doncina_skill_02"
$ws.Range("K18").Value = 1

# --- Highlight the still-blank / carried-over columns (E:K) with the --
# --- purple fill + white text used elsewhere in the sheet for these ---
# --- "unused socket" style rows ----------------------------------------
$ws.Range("E17:K17").Interior.Color = 10498160
$ws.Range("E17:K17").Font.ThemeColor = 2

$ws.Range("E18:K18").Interior.Color = 10498160
$ws.Range("E18:K18").Font.ThemeColor = 2

# --- Restore the selection to where the author left off ----------------
$ws.Range("E16").Select()
